# Update the quarterly database: drop the oldest quarter (1399/06), shift
# all existing quarterly figures one column to the left, and append the
# newly reported quarter (1401/12). Also correct the company name and a
# handful of figures that changed with the new "read_price" algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name correction -------------------------------------------
$ws.Range("B5").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# --- Quarter headers (row 8 and row 24) ---------------------------------
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i   # column E = 5 ... N = 14
    $ws.Cells.Item(8, $col).Value = $quarters[$i]
    $ws.Cells.Item(24, $col).Value = $quarters[$i]
}

# --- Quarterly data rows (shift left by one quarter + newest value) ----
$rowValues = @{
    13 = @(-203, 2564, 3, 0, 2, 6815, 51, 268, 1188, -157)
    14 = @(0, 0, 0, 1310, 1418, -2728, 0, 0, 0, 0)
    15 = @(173, 241, 241, 383, 557, 332, 724, 230, 324, 321)
    16 = @(198, 223, 858, 624, 997, 1012, 1417, 1043, 1590, 2925)
    17 = @(5009, 11226, 15506, 16464, 21467, 22639, 23798, 29616, 21175, 22616)
    19 = @(8445, 24744, 7926, 91333, 24844, 105014, 24418, 79795, 56428, 94667)
    20 = @(13622, 38998, 24534, 110114, 49285, 133084, 50408, 110952, 80705, 120372)
    26 = @(95, 25, 100, 25, 25, 27, 25, 27, 27, 27)
    27 = @(68, 141, 66, 156, 160, 150, 153, 158, 160, 156)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i   # column E = 5 ... N = 14
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
